$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.229.05'
$ws.Range("E2").Value = '  -0.18%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.563.86'
$ws.Range("E3").Value = '  +1.36%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.02'
$ws.Range("E5").Value = '  -0.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.78'
$ws.Range("E6").Value = '  -1.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.562.92'
$ws.Range("E7").Value = '  +1.16%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  +2.37%  '
$ws.Range("E10").Value = '  -0.22%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.82'
$ws.Range("E11").Value = '  -2.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.412'
$ws.Range("E12").Value = '  -0.25%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.167.13'
$ws.Range("E13").Value = '  +1.42%  '
$ws.Range("E14").Value = '  -0.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '30.23'
$ws.Range("E15").Value = '  -0.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.569.51'
$ws.Range("E16").Value = '  +1.67%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.276.15'
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("E18").Value = '  -0.63%  '
$ws.Range("E19").Value = '  +6.50%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.20'
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.73'
$ws.Range("E21").Value = '  -1.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '429.99'
$ws.Range("E22").Value = '  +0.86%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.609'
$ws.Range("E23").Value = '  +1.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.54'
$ws.Range("E24").Value = '  +1.83%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.704.51'
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000117'
$ws.Range("E27").Value = '  -2.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.50'
$ws.Range("E28").Value = '  +0.71%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.08'
$ws.Range("E29").Value = '  -2.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.87'
$ws.Range("E30").Value = '  -1.84%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  -0.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.559.56'
$ws.Range("E32").Value = '  +1.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '25.45'
$ws.Range("E33").Value = '  +0.76%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.45'
$ws.Range("E34").Value = '  -1.84%  '
$ws.Range("E35").Value = '  -8.58%  '
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.82'
$ws.Range("E37").Value = '  +0.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.71'
$ws.Range("E38").Value = '  -1.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.56'
$ws.Range("E39").Value = '  -1.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '173.80'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0849'
$ws.Range("E41").Value = '  -1.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.19'
$ws.Range("E42").Value = '  +0.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.889'
$ws.Range("E43").Value = '  +0.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.93'
$ws.Range("E44").Value = '  +1.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '45.92'
$ws.Range("E45").Value = '  +0.96%  '
$ws.Range("E47").Value = '  -1.85%  '
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.42'
$ws.Range("E48").Value = '  -0.67%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '24.81'
$ws.Range("E49").Value = '  -4.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.12'
$ws.Range("E50").Value = '  -0.87%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.24'
$ws.Range("E51").Value = '  +3.03%  '
